$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "25.877.10"
$ws.Range("E2").Value = "  -0.14%  "
$ws.Range("D3").Value = "1.632.70"
$ws.Range("E3").Value = "  -0.30%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").Value = "'214.15"
$ws.Range("E5").Value = "  -0.46%  "
$ws.Range("D6").Value = "'0.504"
$ws.Range("E6").Value = "  -0.31%  "
$ws.Range("E7").Value = "  +0.01%  "
$ws.Range("D8").Value = "'0.0634"
$ws.Range("E8").Value = "  -0.72%  "
$ws.Range("E9").Value = "  -1.38%  "
$ws.Range("D10").Value = "'19.34"
$ws.Range("E10").Value = "  -1.43%  "
$ws.Range("D11").Value = "'0.0791"
$ws.Range("E11").Value = "  -0.25%  "
$ws.Range("D12").Value = "'4.22"
$ws.Range("E12").Value = "  -0.74%  "
$ws.Range("D13").Value = "1.635.76"
$ws.Range("E13").Value = "  +0.81%  "
$ws.Range("D14").Value = "'0.538"
$ws.Range("E14").Value = "  -1.09%  "
$ws.Range("B15").Value = "ShibaInu"
$ws.Range("C15").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D15").Value = "0.0₃0756"
$ws.Range("E15").Value = "  -1.20%  "
$ws.Range("B16").Value = "Litecoin"
$ws.Range("C16").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D16").Value = "'62.81"
$ws.Range("E16").Value = "  -0.09%  "
$ws.Range("D17").Value = "25.914.21"
$ws.Range("E17").Value = "  -0.04%  "
$ws.Range("E18").Value = "  +0.03%  "
$ws.Range("D19").Value = "'192.73"
$ws.Range("E19").Value = "  -0.23%  "
$ws.Range("D20").Value = "'4.28"
$ws.Range("E20").Value = "  -2.15%  "
$ws.Range("D21").Value = "'9.68"
$ws.Range("E21").Value = "  -2.47%  "
$ws.Range("D22").Value = "'6.12"
$ws.Range("E22").Value = "  -2.58%  "
$ws.Range("E23").Value = "  +3.22%  "
$ws.Range("D24").Value = "'143.56"
$ws.Range("E24").Value = "  -0.30%  "
$ws.Range("D25").Value = "'1.76"
$ws.Range("E25").Value = "  -0.59%  "
$ws.Range("E26").Value = "  -0.07%  "
$ws.Range("D27").Value = "'6.80"
$ws.Range("E27").Value = "  -0.62%  "
$ws.Range("D28").Value = "'15.42"
$ws.Range("E28").Value = "  -0.65%  "
$ws.Range("E29").Value = "  -0.34%  "
$ws.Range("E30").Value = "  -2.51%  "
$ws.Range("E31").Value = "  +0.00%  "
$ws.Range("E32").Value = "  -1.84%  "
$ws.Range("D33").Value = "'1.52"
$ws.Range("E33").Value = "  -1.00%  "
$ws.Range("E34").Value = "  +0.67%  "
$ws.Range("D35").Value = "'0.894"
$ws.Range("E35").Value = "  -0.87%  "
$ws.Range("D36").Value = "1.121.92"
$ws.Range("E36").Value = "  -1.35%  "
$ws.Range("D37").Value = "'0.532"
$ws.Range("E37").Value = "  -2.26%  "
$ws.Range("D38").Value = "'2.45"
$ws.Range("E38").Value = "  -0.41%  "
$ws.Range("D39").Value = "'0.0155"
$ws.Range("E39").Value = "  -1.54%  "
$ws.Range("B40").Value = "Quant"
$ws.Range("C40").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D40").Value = "'98.02"
$ws.Range("E40").Value = "  -1.44%  "
$ws.Range("B41").Value = "TrustWalletToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D41").Value = "'0.792"
$ws.Range("D42").Value = "'5.31"
$ws.Range("E42").Value = "  -3.13%  "
$ws.Range("D43").Value = "0.0₆0114"
$ws.Range("E43").Value = "  -1.07%  "
$ws.Range("D44").Value = "'55.87"
$ws.Range("E44").Value = "  -1.33%  "
$ws.Range("D45").Value = "'1.49"
$ws.Range("E45").Value = "  +1.06%  "
$ws.Range("D46").Value = "'0.0519"
$ws.Range("E46").Value = "  -2.47%  "
$ws.Range("D47").Value = "'7.68"
$ws.Range("E47").Value = "  +0.42%  "
$ws.Range("D48").Value = "'0.412"
$ws.Range("E48").Value = "  -0.55%  "
$ws.Range("E49").Value = "  +0.06%  "
$ws.Range("D50").Value = "'0.0937"
$ws.Range("E50").Value = "  -3.25%  "
$ws.Range("B51").Value = "SynthetixNetwork"
$ws.Range("C51").Value = "https://coinranking.com/coin/sgxZRXbK0FDc+synthetixnetwork-snx"
$ws.Range("D51").Value = "'2.18"
$ws.Range("E51").Value = "  -1.21%  "
